$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-14 -> 2023-09-15) for every data row (rows 2-27).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
